$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239; this shifts the existing rows
# 239-303 down to 240-304 (and extends the used range to R304),
# matching the dimension change A1:R303 -> A1:R304 in the diff.
$ws.Rows(239).Insert()

# Populate the newly inserted row 239 with the new record.
$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44551
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = 100112003
$ws.Cells.Item(239, 7).Value = "Ajo"
$ws.Cells.Item(239, 8).Value = "Chino"
$ws.Cells.Item(239, 9).Value = "1a (cosecha)"
$ws.Cells.Item(239, 10).Value = 85
$ws.Cells.Item(239, 11).Value = 15500
$ws.Cells.Item(239, 12).Value = 16000
$ws.Cells.Item(239, 13).Value = 15765
$ws.Cells.Item(239, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(239, 15).Value = "Limache"
$ws.Cells.Item(239, 16).Value = 1576
$ws.Cells.Item(239, 17).Value = 10
$ws.Cells.Item(239, 18).Value = "Hortaliza"
